$d = $word.ActiveDocument

function Get-PkgXml($innerP) {
    return '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerP + '</w:document></pkg:xmlData></pkg:part></pkg:package>'
}

# Locate the three "G" placeholder paragraphs under the 20/11/2018 heading
# (numId=1) by searching for paragraphs whose visible text is exactly "G".
$gParas = @()
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13)
    if ($t -eq "G") {
        $gParas += $i
    }
}

$idx1 = $gParas[0]
$idx2 = $gParas[1]
$idx3 = $gParas[2]

# --- 1st "G" paragraph -> "Modification de l'interface de fin (score)"
$p1 = $d.Paragraphs.Item($idx1)
$p1.Range.Text = "Modification de l’interface de fin (score)"

# --- 2nd "G" paragraph -> two runs (bookmark added back afterwards)
$p2 = $d.Paragraphs.Item($idx2)
$r2 = $p2.Range
$inner2 = '<w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Maj de la planification (50% </w:t></w:r><w:r><w:t>intégrer)</w:t></w:r></w:p></w:body>'
$r2.InsertXML((Get-PkgXml $inner2))

# --- 3rd "G" paragraph -> text split around a spell-checked "qtCreator"
$p3 = $d.Paragraphs.Item($idx3)
$r3 = $p3.Range
$inner3 = '<w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Recherche sur les différents éléments de </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>qtCreator</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Forms (bouton, image, alignement)</w:t></w:r></w:p></w:body>'
$r3.InsertXML((Get-PkgXml $inner3))

# --- Remove the _GoBack bookmark from its old location (last paragraph of
# the "Difficultés Générales" list) and add text there instead.
$d.Bookmarks.Item("_GoBack").Delete()

# Find the last paragraph in the document (still empty, numId=2 list item)
$lastIdx = $d.Paragraphs.Count
$pLast = $d.Paragraphs.Item($lastIdx)
$rLast = $d.Range($pLast.Range.Start, $pLast.Range.End - 1)
$innerLast = '<w:body><w:p><w:pPr><w:pStyle w:val="Paragraphedeliste"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="2"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Modifier différents éléments sur </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>QtCreator</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Forms</w:t></w:r></w:p></w:body>'
$rLast.InsertXML((Get-PkgXml $innerLast))

# --- Re-add the _GoBack bookmark, collapsed at the end of paragraph 2's
# text (right after "intégrer)"), mirroring its original collapsed position.
$p2 = $d.Paragraphs.Item($idx2)
$r2 = $p2.Range
$bmPos = $r2.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
